$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
$ws.Range("A8").Value = 'Each'
$ws.Range("C8").Value = 57
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.0'
$ws.Range("E8").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F8").Value = 23
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '1311.00'
# --- Row 9 ---
$ws.Range("C9").Value = 11
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '4.0'
$ws.Range("E9").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F9").Value = 50
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '550.00'
# --- Row 10 ---
$ws.Range("C10").Value = 23
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.0'
$ws.Range("E10").Value = 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F10").Value = 78
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '1794.00'
# --- Row 11 ---
$ws.Range("A11").Value = 'Each'
$ws.Range("C11").Value = 13
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '10.0'
$ws.Range("E11").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F11").Value = 303
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '3939.00'
# --- Row 12 ---
$ws.Range("C12").Value = 8
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '11.0'
$ws.Range("E12").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
# --- Row 13 ---
$ws.Range("A13").Value = 'Mtr.'
$ws.Range("C13").Value = 74
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19'
$ws.Range("E13").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F13").Value = 81
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '5994.00'
# --- Row 14 ---
$ws.Range("A14").Value = 'Set'
$ws.Range("C14").Value = 98
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '13.0'
$ws.Range("E14").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F14").Value = 5733
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '561834.00'
# --- Row 15 ---
$ws.Range("C15").Value = 29
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '17.0'
$ws.Range("E15").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
# --- Row 16 ---
$ws.Range("C16").Value = 4
# --- Row 18 ---
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '575422.00'
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = '575422.00'
# --- Row 20 ---
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '575422.00'
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = '575422.00'
